$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I "checked" markers - mirrors RX/TX verification pass noted in the
# commit message ("ok" for items confirmed, explicit quantities where the
# order count itself was double-checked/corrected).
$ws.Range("I2").Value = "ok"
$ws.Range("I3").Value = "ok"
$ws.Range("I4").Value = 60
$ws.Range("I5").Value = "ok"
$ws.Range("I6").Value = "ok"
$ws.Range("I7").Value = "ok"
$ws.Range("I8").Value = "ok"
$ws.Range("I9").Value = "ok"
$ws.Range("I10").Value = "ok"
$ws.Range("I11").Value = "ok"
$ws.Range("I12").Value = "ok"
$ws.Range("I13").Value = "ok"
$ws.Range("I14").Value = "ok"
$ws.Range("I15").Value = "ok"
$ws.Range("I16").Value = "ok"
$ws.Range("I17").Value = "ok"
$ws.Range("I18").Value = "ok"
$ws.Range("I19").Value = 6
$ws.Range("I20").Value = "ok"
$ws.Range("I21").Value = "ok"
$ws.Range("I22").Value = "ok"
$ws.Range("I23").Value = 50
$ws.Range("I24").Value = 250
$ws.Range("I25").Value = 60

# Match the updated view state captured in the saved file: zoomed to 130%
# and the last active selection moved to F23.
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 130
$ws.Range("F23").Select() | Out-Null
